# Update "想去人数" (want-to-go count) figures on both the "展览" and
# "全部类型" worksheets, which carry duplicate data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1889
    $ws.Range("F4").Value = 1151
    $ws.Range("F5").Value = 1178
    $ws.Range("F7").Value = 5981
}
